# Sample-data workbook: split the single "create sample data" helper into
# separate functions. This cosmetic Excel-side trace of that change is:
#   - mark the sender latitude values (B2:B3) with a dedicated decimal
#     number format so they are generated by their own formatting call
#   - widen the (new) third column so the table reads cleanly
#   - leave the cursor on the newly touched cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a 7-decimal numeric format to the sender latitude sample values.
$ws.Range("B2:B3").NumberFormat = "0.0000000"

# Widen column C to fit the recipient/sender longitude labels.
# (ColumnWidth is in characters; the runtime adds a constant padding of
# ~0.8333 characters when serializing to the stored <col width> value, so
# back that out to land on a stored width of 27.)
$ws.Columns.Item(3).ColumnWidth = 26.166666666666668

# Match the cursor position left behind by the edit.
$ws.Range("J15").Select()
